# Build site at 2021-01-29 15:15:03 UTC
# Update course sheet LOQ4217: activation date, objectives, short syllabus,
# full syllabus, evaluation criteria, recovery norm, and bibliography.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Ativação:" date value (kept as literal text, not an Excel date)
$ativacao = '01/01/2021'
$ws.Range("B8").Formula = '=TEXT("' + $ativacao + '","@")'
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4163) | Out-Null
$ws.Range("C8").Formula = '=TEXT("' + $ativacao + '","@")'
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

# Objetivos
$ws.Range("B10").Value = 'Apresentar os conceitos de Logística, Logística Reversa e Gestão da Cadeia de Suprimentos. Capacitar o aluno para aplicação de técnicas e métodos quantitativos para otimização dos problemas em Logística e Cadeias de Suprimentos.'
$ws.Range("C10").Value = 'Apresentar os conceitos de Logística, Logística Reversa e Gestão da Cadeia de Suprimentos. Capacitar o aluno para aplicação de técnicas e métodos quantitativos para otimização dos problemas em Logística e Cadeias de Suprimentos.'

# Programa resumido
$ws.Range("B14").Value = '1. Introdução: 2. Gestão estratégica3. Gestão dos relacionamentos4. Gestão global de suprimentos5. Avaliação de desempenho6. Mapeamento e análise de processos7. Gestão de demanda8. Gestão e coordenação de estoques9. Gestão da logística10. Logística reversa'
$ws.Range("C14").Value = '1. Introdução: 2. Gestão estratégica3. Gestão dos relacionamentos4. Gestão global de suprimentos5. Avaliação de desempenho6. Mapeamento e análise de processos7. Gestão de demanda8. Gestão e coordenação de estoques9. Gestão da logística10. Logística reversa'

# Programa (full syllabus)
$ws.Range("B16").Value = '1. Introdução: A concorrência entre cadeias de suprimento. Definição operacional. A globalização e a gestão de cadeia de suprimentos. Governança das cadeias de suprimentos2. Gestão estratégica: Estratégia de cadeia de suprimentos. Produtos funcionais x produtos inovadores. Fluxos empurrados puxados e híbridos. Custo de transação e a decisão estratégica de comprar ou fazer. Padronização. Integração de parceiros da cadeia de suprimento no projeto de novos produtos e processos.3. Gestão dos relacionamentos: Confiança entre parceiros. Negociação. Gestão do relacionamento com clientes. Segmentação de produtos. Gestão do relacionamento com fornecedores4. Gestão global de suprimentos: Tipos de suplemento. Estrutura organizacional para suprimentos. O processo de suprimento. Coopetição. Ética e responsabilidade social na gestão global de suprimentos5. Avaliação de desempenho: O que é medição de desempenho? Porque medir desempenho. Características de uma boa medida de desempenho. Alinhamento de incentivos em cadeias globais de suprimento. Tipos de contrato de relacionamento6. Mapeamento e análise de processos: Principais processos na cadeia de suprimento. O modelo SCOR (Supply Chain Operations Reference). Análise e melhoramento de processos.7. Gestão de demanda: Ações sobre a demanda para redução de variabilidade. Causas da variabilidade da demanda. Previsão de demanda. Processo de previsão de vendas. Métodos usados em previsões. Método Delphi. Incerteza de previsão8. Gestão e coordenação de estoques: Definição de estoques. Causa do surgimento dos estoques. Tipos de estoque. VMI (vendor management inventory) - estoque gerenciado pelo distribuidor. VOI (vendor owner inventory) - consignação9. Gestão da logística: Centralização versus descentralização. Pontos de armazenagem/distribuição. Funções dos armazéns. Sistemas logísticos escalonados. Localização de unidades logísticas. Gestão de transportes na cadeia de suprimentos.10. Logística reversa: Conceito, importância, estrutura e tendências. Sustentabilidade. Ciclo fechado. Tipos de ciclo fechado. Motivação empresarial. Gerenciamento integrado de resíduos. Modelos de roteirização. Programação de frotas de veículos.'
$ws.Range("C16").Value = '1. Introdução: A concorrência entre cadeias de suprimento. Definição operacional. A globalização e a gestão de cadeia de suprimentos. Governança das cadeias de suprimentos2. Gestão estratégica: Estratégia de cadeia de suprimentos. Produtos funcionais x produtos inovadores. Fluxos empurrados puxados e híbridos. Custo de transação e a decisão estratégica de comprar ou fazer. Padronização. Integração de parceiros da cadeia de suprimento no projeto de novos produtos e processos.3. Gestão dos relacionamentos: Confiança entre parceiros. Negociação. Gestão do relacionamento com clientes. Segmentação de produtos. Gestão do relacionamento com fornecedores4. Gestão global de suprimentos: Tipos de suplemento. Estrutura organizacional para suprimentos. O processo de suprimento. Coopetição. Ética e responsabilidade social na gestão global de suprimentos5. Avaliação de desempenho: O que é medição de desempenho? Porque medir desempenho. Características de uma boa medida de desempenho. Alinhamento de incentivos em cadeias globais de suprimento. Tipos de contrato de relacionamento6. Mapeamento e análise de processos: Principais processos na cadeia de suprimento. O modelo SCOR (Supply Chain Operations Reference). Análise e melhoramento de processos.7. Gestão de demanda: Ações sobre a demanda para redução de variabilidade. Causas da variabilidade da demanda. Previsão de demanda. Processo de previsão de vendas. Métodos usados em previsões. Método Delphi. Incerteza de previsão8. Gestão e coordenação de estoques: Definição de estoques. Causa do surgimento dos estoques. Tipos de estoque. VMI (vendor management inventory) - estoque gerenciado pelo distribuidor. VOI (vendor owner inventory) - consignação9. Gestão da logística: Centralização versus descentralização. Pontos de armazenagem/distribuição. Funções dos armazéns. Sistemas logísticos escalonados. Localização de unidades logísticas. Gestão de transportes na cadeia de suprimentos.10. Logística reversa: Conceito, importância, estrutura e tendências. Sustentabilidade. Ciclo fechado. Tipos de ciclo fechado. Motivação empresarial. Gerenciamento integrado de resíduos. Modelos de roteirização. Programação de frotas de veículos.'

# Critério
$ws.Range("B20").Value = 'Média de Provas e trabalhos (MF).'
$ws.Range("C20").Value = 'Média de Provas e trabalhos (MF).'

# Norma de recuperação
$ws.Range("B21").Value = 'Prova de Recuperação (PR). A Nota final (NF) será a média aritmética entre MF e PR'
$ws.Range("C21").Value = 'Prova de Recuperação (PR). A Nota final (NF) será a média aritmética entre MF e PR'

# Bibliografia
$ws.Range("B22").Value = 'CORRÊA, HENRIQUE LUIZ. Gestão de rede de suprimento: integrando cadeias de suprimento no mundo globalizado. Editora Atlas, 2009CORREA, HENRIQUE LUIZ. Administração de cadeias de suprimento e logística: o essencial. Editora Atlas 2014PIRES, SÉRGIO. Gestão da cadeia de suprimentos (Supply Chain Management): conceitos, estratégias, práticas e casos. Editora Atlas segunda edição. 2009IYER, ANANTH; SESHHADRI, SHIDHAR; VASHER, ROY. A gestão da cadeia de suprimentos da Toyota. Bookman. 2009MYERSON, PAUL. Lean Supply Chain and logistics management. McGrawHill. 2012'
$ws.Range("C22").Value = 'CORRÊA, HENRIQUE LUIZ. Gestão de rede de suprimento: integrando cadeias de suprimento no mundo globalizado. Editora Atlas, 2009CORREA, HENRIQUE LUIZ. Administração de cadeias de suprimento e logística: o essencial. Editora Atlas 2014PIRES, SÉRGIO. Gestão da cadeia de suprimentos (Supply Chain Management): conceitos, estratégias, práticas e casos. Editora Atlas segunda edição. 2009IYER, ANANTH; SESHHADRI, SHIDHAR; VASHER, ROY. A gestão da cadeia de suprimentos da Toyota. Bookman. 2009MYERSON, PAUL. Lean Supply Chain and logistics management. McGrawHill. 2012'
